# Auto-generated cell value updates derived from the authoritative OOXML diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 46613.453
$ws.Range("I6").Value = 63305.875
$ws.Range("J6").Value = 2100.3333
$ws.Range("K6").Value = 189917.625
$ws.Range("L6").Value = 6300.999899999999
$ws.Range("M6").Value = -189805.625
$ws.Range("N6").Value = -6524.999899999999
$ws.Range("H43").Value = 3266.5
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
$ws.Range("H116").Value = 4099.125
$ws.Range("I116").Value = 3565.3333
$ws.Range("K116").Value = 3565.3333
$ws.Range("M116").Value = -123.3332999999998
$ws.Range("H137").Value = 2342.3333
$ws.Range("I137").Value = 2077.2
$ws.Range("J137").Value = 2673.75
$ws.Range("K137").Value = 6231.599999999999
$ws.Range("L137").Value = 8021.25
$ws.Range("M137").Value = -3681.599999999999
$ws.Range("N137").Value = -13121.25
$ws.Range("H138").Value = 2816.75
$ws.Range("I138").Value = 2688.25
$ws.Range("K138").Value = 8064.75
$ws.Range("M138").Value = -2924.75
$ws.Range("H141").Value = 4475.1055
$ws.Range("I141").Value = 4531
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 13593
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -8413
$ws.Range("N141").Value = -22360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H45").Value = 53795.49
$ws.Range("I45").Value = 73502.71000000001
$ws.Range("K45").Value = 73502.71000000001
$ws.Range("M45").Value = -73125.71000000001
$ws.Range("H61").Value = 6951386
$ws.Range("I61").Value = 11116253
$ws.Range("J61").Value = 9940.444
$ws.Range("K61").Value = 11116253
$ws.Range("L61").Value = 9940.444
$ws.Range("M61").Value = -11116041
$ws.Range("N61").Value = -10364.444
$ws.Range("H74").Value = 9642.429
$ws.Range("I74").Value = 7499.4
$ws.Range("K74").Value = 7499.4
$ws.Range("M74").Value = -6625.4
$ws.Range("H77").Value = 9642.429
$ws.Range("I77").Value = 7499.4
$ws.Range("K77").Value = 37497
$ws.Range("M77").Value = -33129
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H122").Value = 3833.3333
$ws.Range("I122").Value = 3591.6667
$ws.Range("K122").Value = 10775.0001
$ws.Range("M122").Value = -8325.000100000001
$ws.Range("H136").Value = 6951386
$ws.Range("I136").Value = 11116253
$ws.Range("J136").Value = 9940.444
$ws.Range("K136").Value = 33348759
$ws.Range("L136").Value = 29821.332
$ws.Range("M136").Value = -33346209
$ws.Range("N136").Value = -34921.33199999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 9282.666999999999
$ws.Range("I36").Value = 9282.666999999999
$ws.Range("K36").Value = 9282.666999999999
$ws.Range("M36").Value = -8748.666999999999
$ws.Range("H37").Value = 1768
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H105").Value = 5810.579
$ws.Range("I105").Value = 5450.5835
$ws.Range("J105").Value = 6427.7144
$ws.Range("K105").Value = 5450.5835
$ws.Range("L105").Value = 6427.7144
$ws.Range("M105").Value = -3703.5835
$ws.Range("N105").Value = -9921.714400000001
$ws.Range("H134").Value = 7049.0625
$ws.Range("I134").Value = 7065.8
$ws.Range("K134").Value = 21197.4
$ws.Range("M134").Value = -18662.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 7498.1665
$ws.Range("I15").Value = 2497.5
$ws.Range("J15").Value = 17499.5
$ws.Range("K15").Value = 2497.5
$ws.Range("L15").Value = 17499.5
$ws.Range("M15").Value = -2327.5
$ws.Range("N15").Value = -17839.5
$ws.Range("H58").Value = 6018.8
$ws.Range("I58").Value = 3319
$ws.Range("J58").Value = 8381.125
$ws.Range("K58").Value = 3319
$ws.Range("L58").Value = 8381.125
$ws.Range("M58").Value = -3116
$ws.Range("N58").Value = -8787.125
$ws.Range("H74").Value = 42027
$ws.Range("J74").Value = 42027
$ws.Range("L74").Value = 42027
$ws.Range("N74").Value = -43775
$ws.Range("H77").Value = 42027
$ws.Range("J77").Value = 42027
$ws.Range("L77").Value = 126081
$ws.Range("N77").Value = -134817
$ws.Range("H86").Value = 6998.5
$ws.Range("I86").Value = 6998.5
$ws.Range("K86").Value = 6998.5
$ws.Range("M86").Value = -5875.5
$ws.Range("H89").Value = 6998.5
$ws.Range("I89").Value = 6998.5
$ws.Range("K89").Value = 34992.5
$ws.Range("M89").Value = -29376.5
$ws.Range("H107").Value = 2024.8235
$ws.Range("I107").Value = 642.7
$ws.Range("K107").Value = 642.7
$ws.Range("M107").Value = 1277.3
$ws.Range("H112").Value = 68589.8
$ws.Range("J112").Value = 68589.8
$ws.Range("L112").Value = 68589.8
$ws.Range("N112").Value = -71543.8
$ws.Range("H134").Value = 7192.778
$ws.Range("I134").Value = 3289.1667
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 9867.500100000001
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -7332.500100000001
$ws.Range("N134").Value = -50070
$ws.Range("H136").Value = 6018.8
$ws.Range("I136").Value = 3319
$ws.Range("J136").Value = 8381.125
$ws.Range("K136").Value = 9957
$ws.Range("L136").Value = 25143.375
$ws.Range("M136").Value = -7407
$ws.Range("N136").Value = -30243.375
$ws.Range("H141").Value = 36288
$ws.Range("J141").Value = 36111.625
$ws.Range("L141").Value = 36111.625
$ws.Range("N141").Value = -46471.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3376.0715
$ws.Range("J5").Value = 3564.3333
$ws.Range("L5").Value = 10692.9999
$ws.Range("N5").Value = -10916.9999
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H92").Value = 640.3333
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 1121
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 3363
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -5859
$ws.Range("H135").Value = 3376.0715
$ws.Range("J135").Value = 3564.3333
$ws.Range("L135").Value = 32078.9997
$ws.Range("N135").Value = -37148.9997
$ws.Range("H137").Value = 12564.5
$ws.Range("J137").Value = 16116
$ws.Range("L137").Value = 48348
$ws.Range("N137").Value = -58548
$ws.Range("H139").Value = 2938.647
$ws.Range("I139").Value = 2400.7273
$ws.Range("J139").Value = 3924.8333
$ws.Range("K139").Value = 7202.1819
$ws.Range("L139").Value = 11774.4999
$ws.Range("M139").Value = -2062.1819
$ws.Range("N139").Value = -22054.4999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 16005.6
$ws.Range("J24").Value = 16005.6
$ws.Range("L24").Value = 16005.6
$ws.Range("N24").Value = -16351.6
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 5595.222
$ws.Range("I132").Value = 5585.3335
$ws.Range("J132").Value = 5615
$ws.Range("K132").Value = 16756.0005
$ws.Range("L132").Value = 16845
$ws.Range("M132").Value = -14226.0005
$ws.Range("N132").Value = -21905
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 668166.7
$ws.Range("I19").Value = 2000
$ws.Range("J19").Value = 1001250
$ws.Range("K19").Value = 2000
$ws.Range("L19").Value = 1001250
$ws.Range("M19").Value = -1830
$ws.Range("N19").Value = -1001590
$ws.Range("H132").Value = 9124.235000000001
$ws.Range("I132").Value = 9274.134
$ws.Range("K132").Value = 27822.402
$ws.Range("M132").Value = -25292.402
$ws.Range("H136").Value = 5939.8
$ws.Range("I136").Value = 4899.6665
$ws.Range("K136").Value = 14698.9995
$ws.Range("M136").Value = -12148.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 11113744
$ws.Range("I2").Value = 20002338
$ws.Range("K2").Value = 20002338
$ws.Range("M2").Value = -20002226
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H107").Value = 1184.8235
$ws.Range("J107").Value = 4175
$ws.Range("L107").Value = 12525
$ws.Range("N107").Value = -16365
$ws.Range("H113").Value = 463.15384
$ws.Range("I113").Value = 463.15384
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1389.46152
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 780.5384799999999
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 130998
$ws.Range("J116").Value = 130998
$ws.Range("L116").Value = 130998
$ws.Range("N116").Value = -140176
$ws.Range("H132").Value = 3826.131
$ws.Range("I132").Value = 3389.4255
$ws.Range("K132").Value = 10168.2765
$ws.Range("M132").Value = -7638.2765
$ws.Range("H136").Value = 6219.815
$ws.Range("I136").Value = 5130.8423
$ws.Range("K136").Value = 15392.5269
$ws.Range("M136").Value = -12842.5269
